# Restore C10 on the "Rules" sheet from 18 to 1.
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Rules") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Range("C10").Value = 1
